$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column, matching the style of the other headers
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Save = 1 if sum (column G) > 10, else 0, for rows 2..48
for ($r = 2; $r -le 48; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -gt 10) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
